# The commit swaps the raw contents of ppt/theme/theme1.xml (the "Office
# Theme" used by the notes master) and ppt/theme/theme2.xml (the "Integral"
# theme used by the slide master / all slides): after the edit theme2.xml
# (the theme that actually paints the deck) carries the plain "Office
# Theme" color palette instead of "Integral".
#
# The color palette is the only part that actually differs between the two
# theme parts (font scheme / format scheme are byte-identical), so we
# reproduce the visible effect by rewriting the live theme's 12 scheme
# colors through the ThemeColorScheme object (MsoThemeColorSchemeIndex
# order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) to the values the
# "Office Theme" previously used.

$p = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Colors(1).RGB  = 0        # dk1      000000
$tcs.Colors(2).RGB  = 16777215 # lt1      FFFFFF
$tcs.Colors(3).RGB  = 6968388  # dk2      44546A
$tcs.Colors(4).RGB  = 15132391 # lt2      E7E6E6
$tcs.Colors(5).RGB  = 13998939 # accent1  5B9BD5
$tcs.Colors(6).RGB  = 3243501  # accent2  ED7D31
$tcs.Colors(7).RGB  = 10855845 # accent3  A5A5A5
$tcs.Colors(8).RGB  = 49407    # accent4  FFC000
$tcs.Colors(9).RGB  = 12874308 # accent5  4472C4
$tcs.Colors(10).RGB = 4697456  # accent6  70AD47
$tcs.Colors(11).RGB = 12673797 # hlink    0563C1
$tcs.Colors(12).RGB = 7491477  # folHlink 954F72
